$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B-column values first so shared-string indices 2,3,4 are allocated
# in the order: dummy, kaggle_bm, kaggle param my dataset
$ws.Range("B2").Value = "dummy"
$ws.Range("B3").Value = "kaggle_bm"
$ws.Range("B4").Value = "kaggle param my dataset"

# Then set C-column values so shared-string indices continue:
# dummy (reuses index 2), then new indices 5,6 for the two multi-line strings
$ws.Range("C2").Value = "dummy"
$ws.Range("C3").Value = "`n    kaggle params`n    "
$ws.Range("C4").Value = "`n    kaggle params with my dataset`n    "

# Remove row 5 entirely (was A5=3, B5="linear regression", C5="1 row lookback")
$ws.Rows.Item(5).Delete()
